$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

$ws.Range("A122").Value = "125"
$ws.Range("B122").Value = 501
$ws.Range("C122").Value = 322
$ws.Range("D122").Value = 1349
$ws.Range("E122").Value = 370
$ws.Range("F122").Value = "you have already watched this story dialog -> title"

$ws.Range("A123").Value = "126"
$ws.Range("B123").Value = 631
$ws.Range("C123").Value = 680
$ws.Range("D123").Value = 809
$ws.Range("E123").Value = 728
$ws.Range("F123").Value = "you have already watched this story dialog -> cancel button"

# A column keeps the centered text-format style (same as A121), F column
# keeps the centered general-format style (same as F121); B:E are left
# unstyled (default "Normal" style), matching the author's edit (only A
# and F carry explicit styles on the new rows).
$aRange = $ws.Range("A122:A123")
$aRange.HorizontalAlignment = -4108
$aRange.VerticalAlignment = -4108
$aRange.NumberFormat = "@"

$fRange = $ws.Range("F122:F123")
$fRange.HorizontalAlignment = -4108
$fRange.VerticalAlignment = -4108

$ws.Range("B122:E123").Style = "Normal"

$ws.Range("F124").Select()
